$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 26 de Marzo de 2020 a las 02:12"

# Apply updated COVID-19 counts; table is kept sorted by Casos totales (col B) desc,
# so several rows pick up a new Pais (col A) along with refreshed stats.
# Row 6: Estados Unidos
$ws.Cells.Item(6,2).Value = 65998
$ws.Cells.Item(6,3).Value = 11142
$ws.Cells.Item(6,5).Value = 64661
$ws.Cells.Item(6,6).Value = 1452
$ws.Cells.Item(6,7).Value = 163
$ws.Cells.Item(6,8).Value = 943
# Row 117: Consejo Danes para los Refugiados
$ws.Cells.Item(117,1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(117,3).Value = 3
$ws.Cells.Item(117,6).Value = 0
# Row 118: Mauricio
$ws.Cells.Item(118,1).Value = "Mauricio"
$ws.Cells.Item(118,3).Value = 6
$ws.Cells.Item(118,6).Value = 1
# Row 123: Bolivia
$ws.Cells.Item(123,1).Value = "Bolivia"
$ws.Cells.Item(123,2).Value = 38
$ws.Cells.Item(123,3).Value = 9
$ws.Cells.Item(123,5).Value = 38
$ws.Cells.Item(123,6).Value = 0
$ws.Cells.Item(123,7).Value = 0
$ws.Cells.Item(123,8).Value = 0
# Row 124: Paraguay
$ws.Cells.Item(124,1).Value = "Paraguay"
$ws.Cells.Item(124,2).Value = 37
$ws.Cells.Item(124,3).Value = 10
$ws.Cells.Item(124,5).Value = 34
$ws.Cells.Item(124,6).Value = 1
$ws.Cells.Item(124,7).Value = 1
$ws.Cells.Item(124,8).Value = 3
# Row 125: Mayotte
$ws.Cells.Item(125,1).Value = "Mayotte"
$ws.Cells.Item(125,2).Value = 36
$ws.Cells.Item(125,3).Value = 0
$ws.Cells.Item(125,5).Value = 36
# Row 141: Nueva Caledonia
$ws.Cells.Item(141,1).Value = "Nueva Caledonia"
$ws.Cells.Item(141,3).Value = 4
# Row 142: Uganda
$ws.Cells.Item(142,1).Value = "Uganda"
$ws.Cells.Item(142,3).Value = 5
# Row 144: El Salvador
$ws.Cells.Item(144,1).Value = "El Salvador"
$ws.Cells.Item(144,3).Value = 8
$ws.Cells.Item(144,4).Value = 0
$ws.Cells.Item(144,5).Value = 13
# Row 145: Maldivas
$ws.Cells.Item(145,1).Value = "Maldivas"
$ws.Cells.Item(145,2).Value = 13
$ws.Cells.Item(145,4).Value = 8
$ws.Cells.Item(145,5).Value = 5
# Row 147: Etiopia
$ws.Cells.Item(147,1).Value = "Etiopia"
$ws.Cells.Item(147,2).Value = 12
$ws.Cells.Item(147,3).Value = 0
$ws.Cells.Item(147,5).Value = 12
# Row 148: Republica de Yibuti
$ws.Cells.Item(148,1).Value = "Republica de Yibuti"
$ws.Cells.Item(148,3).Value = 8
# Row 149: San Martin (Parte Francesa)
$ws.Cells.Item(149,1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(149,3).Value = 3
# Row 150: Dominica
$ws.Cells.Item(150,1).Value = "Dominica"
$ws.Cells.Item(150,2).Value = 11
$ws.Cells.Item(150,3).Value = 4
$ws.Cells.Item(150,5).Value = 11
# Row 151: Mongolia
$ws.Cells.Item(151,1).Value = "Mongolia"
$ws.Cells.Item(151,2).Value = 10
$ws.Cells.Item(151,3).Value = 0
$ws.Cells.Item(151,5).Value = 10
# Row 153: Surinam
$ws.Cells.Item(153,1).Value = "Surinam"
# Row 154: Haiti
$ws.Cells.Item(154,1).Value = "Haiti"
# Row 156: Bermudas
$ws.Cells.Item(156,1).Value = "Bermudas"
$ws.Cells.Item(156,3).Value = 1
# Row 157: Seychelles
$ws.Cells.Item(157,1).Value = "Seychelles"
$ws.Cells.Item(157,3).Value = 0
# Row 165: Mozambique
$ws.Cells.Item(165,1).Value = "Mozambique"
$ws.Cells.Item(165,3).Value = 2
# Row 166: Fiyi
$ws.Cells.Item(166,1).Value = "Fiyi"
$ws.Cells.Item(166,3).Value = 1
# Row 167: Bahamas
$ws.Cells.Item(167,1).Value = "Bahamas"
$ws.Cells.Item(167,4).Value = 1
$ws.Cells.Item(167,8).Value = 0
# Row 168: Guyana
$ws.Cells.Item(168,1).Value = "Guyana"
$ws.Cells.Item(168,4).Value = 0
$ws.Cells.Item(168,8).Value = 1
# Row 170: Santa Sede
$ws.Cells.Item(170,1).Value = "Santa Sede"
$ws.Cells.Item(170,3).Value = 0
# Row 171: Guinea
$ws.Cells.Item(171,1).Value = "Guinea"
# Row 172: Congo
$ws.Cells.Item(172,1).Value = "Congo"
# Row 173: Eritrea
$ws.Cells.Item(173,1).Value = "Eritrea"
$ws.Cells.Item(173,3).Value = 3
# Row 175: Angola
$ws.Cells.Item(175,1).Value = "Angola"
# Row 176: Antigua y Barbuda
$ws.Cells.Item(176,1).Value = "Antigua y Barbuda"
# Row 177: San Martin (Parte Holandesa)
$ws.Cells.Item(177,1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(177,3).Value = 1
# Row 178: Birmania
$ws.Cells.Item(178,1).Value = "Birmania"
# Row 179: Republica de Africa Central
$ws.Cells.Item(179,1).Value = "Republica de Africa Central"
$ws.Cells.Item(179,3).Value = 0
# Row 180: Republica del Chad
$ws.Cells.Item(180,1).Value = "Republica del Chad"
$ws.Cells.Item(180,3).Value = 0
# Row 181: Liberia
$ws.Cells.Item(181,1).Value = "Liberia"
# Row 182: Laos
$ws.Cells.Item(182,1).Value = "Laos"
$ws.Cells.Item(182,3).Value = 1
# Row 183: San Bartolome
$ws.Cells.Item(183,1).Value = "San Bartolome"
# Row 185: Sudan
$ws.Cells.Item(185,1).Value = "Sudan"
# Row 186: Gambia
$ws.Cells.Item(186,1).Value = "Gambia"
# Row 187: Nepal
$ws.Cells.Item(187,1).Value = "Nepal"
$ws.Cells.Item(187,3).Value = 1
$ws.Cells.Item(187,4).Value = 1
$ws.Cells.Item(187,8).Value = 0
# Row 188: Zimbabue
$ws.Cells.Item(188,1).Value = "Zimbabue"
$ws.Cells.Item(188,3).Value = 0
$ws.Cells.Item(188,4).Value = 0
$ws.Cells.Item(188,8).Value = 1
# Row 189: Guinea-Bisau
$ws.Cells.Item(189,1).Value = "Guinea-Bisau"
$ws.Cells.Item(189,3).Value = 2
# Row 190: Islas Virgenes Britanicas
$ws.Cells.Item(190,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(190,3).Value = 2
# Row 191: San Cristobal y Nieves
$ws.Cells.Item(191,1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(191,3).Value = 2
# Row 192: Mauritania
$ws.Cells.Item(192,1).Value = "Mauritania"
$ws.Cells.Item(192,3).Value = 0
# Row 194: Belice
$ws.Cells.Item(194,1).Value = "Belice"
$ws.Cells.Item(194,3).Value = 1
# Row 195: Mali
$ws.Cells.Item(195,1).Value = "Mali"
# Row 196: Butan
$ws.Cells.Item(196,1).Value = "Butan"
$ws.Cells.Item(196,3).Value = 0
# Row 197: San Vicente y las Granadinas
$ws.Cells.Item(197,1).Value = "San Vicente y las Granadinas"
# Row 198: Granada
$ws.Cells.Item(198,1).Value = "Granada"
# Row 200: Islas Turcas y Caicos
$ws.Cells.Item(200,1).Value = "Islas Turcas y Caicos"
# Row 201: Somalia
$ws.Cells.Item(201,1).Value = "Somalia"
# Row 202: Papua Nueva Guinea
$ws.Cells.Item(202,1).Value = "Papua Nueva Guinea"
# Row 203: Montserrat
$ws.Cells.Item(203,1).Value = "Montserrat"
# Row 204: Libia
$ws.Cells.Item(204,1).Value = "Libia"
